
# Auto-update draw results on excel 2025-12-12T17:41:58Z
# Adds the new Pick 3 draw row (row 87) to the Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# Force the new row to be stored as plain text (matching every other row in
# the sheet, where dates/phase codes/results are all text, not numbers or
# dates) by setting the number format to "@" before assigning values.
$target = $ws.Range("A" + $row + ":E" + $row)
$target.NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-12-12"
$ws.Range("B" + $row).Value = "Pick 3"
$ws.Range("C" + $row).Value = "251212"
$ws.Range("D" + $row).Value = "8-2-5"
$ws.Range("E" + $row).Value = "2025-12-12T21:41:58.074+04:00"
